$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = " 900 - 1000"
$ws.Range("A7").Select()
